# "add libs and new apk"
#
# The underlying sheet is a price/reference list ("AgroComplit/наименования.xlsx").
# This commit clears out a small scratch to-do list that had been jotted down in
# column D (rows 15-22) of "Лист1", leaving those cells blank, and moves the
# active selection to G22 (where the user's cursor ended up after the cleanup).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear the stray "to do" notes that lived in D15:D20 -------------------
# Those cells used a bold/shaded style (fontId 2, fillId theme9) that is not
# used anywhere else in the sheet. Re-apply the plain "vertical-center" style
# that already exists on the neighbouring E column (style index 5) before
# wiping the text, so the now-empty cells pick up the same formatting as any
# other blank cell in the table instead of dragging along the bespoke fill.
$fmtSource = $ws.Range("E15")
$fmtSource.Copy()
$ws.Range("D15:D20").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("D15").ClearContents()
$ws.Range("D16").ClearContents()
$ws.Range("D17").ClearContents()
$ws.Range("D18").ClearContents()
$ws.Range("D19").ClearContents()
$ws.Range("D20").ClearContents()

# D21/D22 already use the plain bold style (fontId 2, no fill) that is still
# used elsewhere, so only their text needs to go.
$ws.Range("D21").ClearContents()
$ws.Range("D22").ClearContents()

# --- Move the saved selection, matching where editing left the cursor ------
$ws.Range("G22").Select()
